$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.783.90'
$ws.Range('E2').Value = '  -1.45%  '
$ws.Range('D3').Value = '1.547.60'
$ws.Range('E3').Value = '  -1.70%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.54%  '
$ws.Range('E6').Value = '  -1.87%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -1.24%  '
$ws.Range('E9').Value = '  -4.10%  '
$ws.Range('E10').Value = '  -1.43%  '
$ws.Range('E11').Value = '  -1.62%  '
$ws.Range('D12').Value = '1.767.33'
$ws.Range('E12').Value = '  -1.67%  '
$ws.Range('D13').Value = '1.555.57'
$ws.Range('E13').Value = '  -1.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.67'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.511'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.43%  '
$ws.Range('D16').Value = '26.767.01'
$ws.Range('E16').Value = '  -1.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.06'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.96%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '213.35'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('E20').Value = '  -1.63%  '
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('E22').Value = '  -1.19%  '
$ws.Range('E23').Value = '  -5.31%  '
$ws.Range('E24').Value = '  -2.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.18'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('E26').Value = '  -2.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E29').Value = '  -1.48%  '
$ws.Range('E30').Value = '  -0.62%  '
$ws.Range('E31').Value = '  -1.27%  '
$ws.Range('E32').Value = '  +0.27%  '
$ws.Range('D33').Value = '1.345.90'
$ws.Range('E33').Value = '  -4.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.91'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.51'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.29'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.926'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.60%  '
$ws.Range('E38').Value = '  -1.37%  '
$ws.Range('E39').Value = '  +0.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.800'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.90%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.69'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.75%  '
$ws.Range('E43').Value = '  -1.05%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('E45').Value = '  -5.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '62.83'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.45%  '
$ws.Range('D47').Value = '1.682.24'
$ws.Range('E47').Value = '  -1.57%  '
$ws.Range('E48').Value = '  -3.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '85.61'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0515'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.27%  '
$ws.Range('E51').Value = '  -1.90%  '
